$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Copy the existing date formatting (style index for short-date m/d/yyyy) from row 30/31
# down to the new date cells, so a duplicate number format isn't created.
$ws.Range("C30:D31").Copy() | Out-Null
$ws.Range("C32:D33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 32: Deep Work by Cal Newport (G32 filled in later, after row 33)
$ws.Range("A32").Value = "Deep Work"
$ws.Range("B32").Value = "Cal Newport"
$ws.Range("C32").Value = 43888
$ws.Range("D32").Value = 43889
$ws.Range("E32").Value = "self improvement;productivity;excellence;focus"
$ws.Range("F32").Value = "Audio"

# Row 33: Crazy Rich by Jerry Oppenheimer
$ws.Range("A33").Value = "Crazy Rich"
$ws.Range("B33").Value = "Jerry Oppenheimer"
$ws.Range("C33").Value = 43889
$ws.Range("D33").Value = 43892
$ws.Range("E33").Value = "dynasty;johnson & johnson;heirs;scandal;history"
$ws.Range("F33").Value = "Audio"
$ws.Range("G33").Value = "15 Hours 13 Mins"

# Now go back and fill in G32
$ws.Range("G32").Value = "7 Hours 44 Mins"

$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("G33").Select()
